$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$c = $ws1.Range("E23")
$c.ClearFormats()
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.ThemeColor = 1
$c.NumberFormat = "#,##0.00"
$c.HorizontalAlignment = -4152
$c.VerticalAlignment = -4108
$c.WrapText = $True
